$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 174, shifting existing rows 174-271 down to 176-273.
$ws.Range("A174:A175").EntireRow.Insert()

# Populate the two newly inserted rows (174 and 175) with their new data.
$ws.Cells.Item(174, 1).Value = 3
$ws.Cells.Item(174, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(174, 3).Value = 'Coquimbo'
$ws.Cells.Item(174, 4).Value = 44518
$ws.Cells.Item(174, 5).Value = 5
$ws.Cells.Item(174, 6).Value = 100112017
$ws.Cells.Item(174, 7).Value = 'Apio'
$ws.Cells.Item(174, 8).Value = 'Americana (o)'
$ws.Cells.Item(174, 9).Value = 'Primera'
$ws.Cells.Item(174, 10).Value = 160
$ws.Cells.Item(174, 11).Value = 9000
$ws.Cells.Item(174, 12).Value = 9000
$ws.Cells.Item(174, 13).Value = 9000
$ws.Cells.Item(174, 14).Value = '$/docena de matas'
$ws.Cells.Item(174, 15).Value = 'Pan de Azúcar'
$ws.Cells.Item(174, 16).Value = 1500
$ws.Cells.Item(174, 17).Value = 6
$ws.Cells.Item(174, 18).Value = 'Hortaliza'
$ws.Cells.Item(175, 1).Value = 3
$ws.Cells.Item(175, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(175, 3).Value = 'Coquimbo'
$ws.Cells.Item(175, 4).Value = 44518
$ws.Cells.Item(175, 5).Value = 5
$ws.Cells.Item(175, 6).Value = 100112017
$ws.Cells.Item(175, 7).Value = 'Apio'
$ws.Cells.Item(175, 8).Value = 'Americana (o)'
$ws.Cells.Item(175, 9).Value = 'Segunda'
$ws.Cells.Item(175, 10).Value = 60
$ws.Cells.Item(175, 11).Value = 7000
$ws.Cells.Item(175, 12).Value = 7000
$ws.Cells.Item(175, 13).Value = 7000
$ws.Cells.Item(175, 14).Value = '$/docena de matas'
$ws.Cells.Item(175, 15).Value = 'Pan de Azúcar'
$ws.Cells.Item(175, 16).Value = 1167
$ws.Cells.Item(175, 17).Value = 6
$ws.Cells.Item(175, 18).Value = 'Hortaliza'

Write-Host "Applied edits to rows 174-175 (inserted) and verified downstream shift for rows 176-273."
